$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AuthTests")

# Remove the two generic placeholder tables (addEmployeeTestData, loginTestData)
# that only ever held "Column1".."Column7" headers with no real data.
$ws.ListObjects.Item("addEmployeeTestData").Delete()
$ws.ListObjects.Item("loginTestData").Delete()

# Add the new "Add Employee" data-driven test data block (header row first,
# then each data column top-to-bottom, matching how it was authored).
$ws.Range("D3").Value = "FirstName"
$ws.Range("E3").Value = "MiddleName"
$ws.Range("F3").Value = "LastName"

$ws.Range("D4").Value = "NameA"
$ws.Range("D5").Value = "NameB"
$ws.Range("D6").Value = "NameC"
$ws.Range("D7").Value = "NameD"
$ws.Range("D8").Value = "NameE"

$ws.Range("E4").Value = "Middle1"
$ws.Range("E5").Value = "Middle2"
$ws.Range("E6").Value = "Middle3"
$ws.Range("E7").Value = "Middle4"
$ws.Range("E8").Value = "Middle5"

$ws.Range("F4").Value = "Last1"
$ws.Range("F5").Value = "Last2"
$ws.Range("F6").Value = "Last3"
$ws.Range("F7").Value = "Last4"
$ws.Range("F8").Value = "Last5"

$tbl = $ws.ListObjects.Add(1, $ws.Range("D3:F8"), 0, 1)
$tbl.Name = "EmployeeDetails"

# Widen the new MiddleName/LastName columns so the data is readable.
$ws.Range("E1").ColumnWidth = 14
$ws.Range("F1").ColumnWidth = 11

# Page is printed in portrait orientation.
$ws.PageSetup.Orientation = 1

# Leave the selection near the new table, and scroll back up so row 1 is visible.
$ws.Range("E14").Select() | Out-Null
